# feat: add 2022-Q1 data
#
# Before: sheets are 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计(summary)
# After:  sheets are 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计(summary)
#
# The previously-last "总计" sheet becomes the new "2022-Q1" fund-holdings
# sheet (keeping its original sheetId), and a brand new "总计" sheet is
# appended at the end with an extra row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" data sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Template sheet that already has the full 8-column (A:H) layout/styles we
# need (bold+bordered header row B1:H1, bold+bordered index column A).
$template = $wb.Worksheets.Item("2021-Q4")

# Extend the header row style (B1 already carries it) across E1:H1.
$template.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Extend the index-column style (A2 already carries it) down to A7:A17.
$template.Range("A2").Copy()
$q1.Range("A7:A17").PasteSpecial(-4122)

# Wipe the old 2021-Q4/.../2020-Q4 rollup values (B2:D6) before refilling.
$q1.Range("B2:D6").ClearContents()

$fundRows = @(
    @('0','004616','中欧电子信息产业沪港深股票A','14.54','92.26','5.37','0.7808','4'),
    @('1','005763','中欧电子信息产业沪港深股票C','7.73','92.26','5.37','0.4151','4'),
    @('2','000219','博时裕益灵活配置混合','3.90','91.22','5.17','0.2016','2'),
    @('3','000586','景顺长城中小板创业板精选股票','2.42','94.15','5.96','0.1442','8'),
    @('4','013393','信达澳银价值精选混合A','3.61','81.31','3.10','0.1119','5'),
    @('5','001415','信诚新锐回报灵活配置混合A','9.07','24.72','0.69','0.0626','6'),
    @('6','003234','信诚至利灵活配置混合A','8.99','22.05','0.61','0.0548','6'),
    @('7','001402','信诚新选回报灵活配置混合A','8.37','22.05','0.63','0.0527','5'),
    @('8','004157','信诚至诚灵活配置混合A','7.32','22.71','0.63','0.0461','6'),
    @('9','002046','信诚新锐回报灵活配置混合B','5.19','24.72','0.69','0.0358','6'),
    @('10','003235','信诚至利灵活配置混合C','5.30','22.05','0.61','0.0323','6'),
    @('11','002030','信诚新选回报灵活配置混合B','3.85','22.05','0.63','0.0243','5'),
    @('12','004194','招商中证1000指数增强A','1.76','94.40','1.08','0.0190','5'),
    @('13','004158','信诚至诚灵活配置混合B','2.18','22.71','0.63','0.0137','6'),
    @('14','013394','信达澳银价值精选混合C','0.37','81.31','3.10','0.0115','5'),
    @('15','004195','招商中证1000指数增强C','0.68','94.40','1.08','0.0073','5')
)

$r = 2
foreach ($row in $fundRows) {
    # A: numeric row index
    $q1.Cells.Item($r, 1).Value = [double]$row[0]

    # B: fund code, C: fund name -- plain text
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]

    # D..G: numeric-looking figures stored as literal text (matches source)
    for ($col = 4; $col -le 7; $col++) {
        $cell = $q1.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
        $cell.Style = "Normal"
    }

    # H: position-rank, numeric
    $q1.Cells.Item($r, 8).Value = [double]$row[7]

    $r++
}

# ---------------------------------------------------------------------------
# Step 2: append a brand new "总计" sheet (gets a fresh sheetId) with the
# updated rollup table (2022-Q1 plus the previous quarters).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$total.Name = "总计"

# Borrow header/index styling from the 2022-Q1 sheet (same "s=2" style).
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summaryRows = @(
    @('0','2022-Q1','16','2.01'),
    @('1','2021-Q4','25','8.43'),
    @('2','2021-Q3','11','3.59'),
    @('3','2021-Q2','3','0.58'),
    @('4','2021-Q1','16','7.91'),
    @('5','2020-Q4','9','3.12')
)

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = [double]$row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = [double]$row[2]
    $total.Cells.Item($r, 4).Value = [double]$row[3]
    $r++
}

# Restore the originally-active sheet/selection.
$wb.Worksheets.Item(1).Activate()
